# Add team record (Wins/Losses/Ties) columns AD:AF to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered style) by copying the format from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2:40 all get the same team record values.
$lastRow = 40
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
